# NIT-8300264537.xlsx - "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The underlying database value for "Salario Basico" (column F) was corrected:
# the figures for the first ("2107") and last ("2201") "Periodo Mora" rows were
# transposed - swap them back.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$f16 = $ws.Range("F16").Value2
$f22 = $ws.Range("F22").Value2

$ws.Range("F16").Value2 = $f22
$ws.Range("F22").Value2 = $f16
